$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header A1 from "ArticleName" to "ArticleID"
$ws.Range("A1").Value = "ArticleID"

# Remove the obsolete columns L:Q (original/replication r/d/n collected headers)
$ws.Range("L1:Q1").ClearContents()

# Rename study labels in column A (short codes instead of descriptive names)
# Order matches how the new shared strings were introduced upstream
$ws.Range("A7").Value = "Econ"
$ws.Range("A6").Value = "NatSci"
$ws.Range("A3").Value = "ManyLabs1"
$ws.Range("A4").Value = "ManyLabs2"
$ws.Range("A5").Value = "ManyLabs3"
$ws.Range("A8").Value = "xPhi"

# Mark rows 6, 7, 8 as having Data_Cleaned = 1 (column K)
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 1
$ws.Range("K8").Value = 1

# Update the active selection to match the source workbook
$ws.Range("L7").Select()

$wb.Save()
